# Auto-generated script applying numeric corrections to the Leve profit
# tables across all eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below were recomputed by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value2 = 3254.3333  # H70: 3404.4546 -> 3254.3333
$ws.Cells.Item(70, 9).Value2 = 4908.3335  # I70: 5590 -> 4908.3335
$ws.Cells.Item(70, 10).Value2 = 1600.3334  # J70: 1583.1666 -> 1600.3334
$ws.Cells.Item(70, 11).Value2 = 14725.0005  # K70: 16770 -> 14725.0005
$ws.Cells.Item(70, 12).Value2 = 4801.0002  # L70: 4749.4998 -> 4801.0002
$ws.Cells.Item(70, 13).Value2 = -14455.0005  # M70: -16500 -> -14455.0005
$ws.Cells.Item(70, 14).Value2 = -5341.0002  # N70: -5289.4998 -> -5341.0002
$ws.Cells.Item(73, 8).Value2 = 3254.3333  # H73: 3404.4546 -> 3254.3333
$ws.Cells.Item(73, 9).Value2 = 4908.3335  # I73: 5590 -> 4908.3335
$ws.Cells.Item(73, 10).Value2 = 1600.3334  # J73: 1583.1666 -> 1600.3334
$ws.Cells.Item(73, 11).Value2 = 14725.0005  # K73: 16770 -> 14725.0005
$ws.Cells.Item(73, 12).Value2 = 4801.0002  # L73: 4749.4998 -> 4801.0002
$ws.Cells.Item(73, 13).Value2 = -13789.0005  # M73: -15834 -> -13789.0005
$ws.Cells.Item(73, 14).Value2 = -6673.0002  # N73: -6621.4998 -> -6673.0002
$ws.Cells.Item(74, 8).Value2 = 3467.8667  # H74: 3460.8667 -> 3467.8667
$ws.Cells.Item(74, 9).Value2 = 2772.923  # I74: 2756.7693 -> 2772.923
$ws.Cells.Item(74, 11).Value2 = 2772.923  # K74: 2756.7693 -> 2772.923
$ws.Cells.Item(74, 13).Value2 = -1836.923  # M74: -1820.7693 -> -1836.923
$ws.Cells.Item(77, 8).Value2 = 3467.8667  # H77: 3460.8667 -> 3467.8667
$ws.Cells.Item(77, 9).Value2 = 2772.923  # I77: 2756.7693 -> 2772.923
$ws.Cells.Item(77, 11).Value2 = 13864.615  # K77: 13783.8465 -> 13864.615
$ws.Cells.Item(77, 13).Value2 = -9184.614999999998  # M77: -9103.8465 -> -9184.614999999998
$ws.Cells.Item(132, 8).Value2 = 3177.2173  # H132: 2907.5354 -> 3177.2173
$ws.Cells.Item(132, 9).Value2 = 2926.027  # I132: 2618.1235 -> 2926.027
$ws.Cells.Item(132, 11).Value2 = 8778.081  # K132: 7854.370500000001 -> 8778.081
$ws.Cells.Item(132, 13).Value2 = -6248.081  # M132: -5324.370500000001 -> -6248.081
$ws.Cells.Item(137, 8).Value2 = 2853.4092  # H137: 2941.2769 -> 2853.4092
$ws.Cells.Item(137, 9).Value2 = 864.619  # I137: 1190.7368 -> 864.619
$ws.Cells.Item(137, 10).Value2 = 6333.7915  # J137: 5405 -> 6333.7915
$ws.Cells.Item(137, 11).Value2 = 2593.857  # K137: 3572.2104 -> 2593.857
$ws.Cells.Item(137, 12).Value2 = 19001.3745  # L137: 16215 -> 19001.3745
$ws.Cells.Item(137, 13).Value2 = -43.85699999999997  # M137: -1022.2104 -> -43.85699999999997
$ws.Cells.Item(137, 14).Value2 = -24101.3745  # N137: -21315 -> -24101.3745
$ws.Cells.Item(138, 8).Value2 = 1990.9531  # H138: 2058.1875 -> 1990.9531
$ws.Cells.Item(138, 9).Value2 = 806.7692  # I138: 871.2778 -> 806.7692
$ws.Cells.Item(138, 10).Value2 = 3838.28  # J138: 3584.2144 -> 3838.28
$ws.Cells.Item(138, 11).Value2 = 2420.3076  # K138: 2613.8334 -> 2420.3076
$ws.Cells.Item(138, 12).Value2 = 11514.84  # L138: 10752.6432 -> 11514.84
$ws.Cells.Item(138, 13).Value2 = 2719.6924  # M138: 2526.1666 -> 2719.6924
$ws.Cells.Item(138, 14).Value2 = -21794.84  # N138: -21032.6432 -> -21794.84
$ws.Cells.Item(141, 8).Value2 = 3768.0454  # H141: 3567.8723 -> 3768.0454
$ws.Cells.Item(141, 9).Value2 = 1242.5428  # I141: 1242.7142 -> 1242.5428
$ws.Cells.Item(141, 10).Value2 = 13589.444  # J141: 10349.583 -> 13589.444
$ws.Cells.Item(141, 11).Value2 = 3727.6284  # K141: 3728.1426 -> 3727.6284
$ws.Cells.Item(141, 12).Value2 = 40768.33199999999  # L141: 31048.749 -> 40768.33199999999
$ws.Cells.Item(141, 13).Value2 = 1452.3716  # M141: 1451.8574 -> 1452.3716
$ws.Cells.Item(141, 14).Value2 = -51128.33199999999  # N141: -41408.749 -> -51128.33199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value2 = 1145.4286  # H61: 1149.7206 -> 1145.4286
$ws.Cells.Item(61, 9).Value2 = 853.2406999999999  # I61: 865.5848999999999 -> 853.2406999999999
$ws.Cells.Item(61, 10).Value2 = 2131.5625  # J61: 2153.6667 -> 2131.5625
$ws.Cells.Item(61, 11).Value2 = 853.2406999999999  # K61: 865.5848999999999 -> 853.2406999999999
$ws.Cells.Item(61, 12).Value2 = 2131.5625  # L61: 2153.6667 -> 2131.5625
$ws.Cells.Item(61, 13).Value2 = -641.2406999999999  # M61: -653.5848999999999 -> -641.2406999999999
$ws.Cells.Item(61, 14).Value2 = -2555.5625  # N61: -2577.6667 -> -2555.5625
$ws.Cells.Item(74, 8).Value2 = 3717.2778  # H74: 3028.9556 -> 3717.2778
$ws.Cells.Item(74, 9).Value2 = 872.2069  # I74: 714.025 -> 872.2069
$ws.Cells.Item(74, 10).Value2 = 15504  # J74: 21548.4 -> 15504
$ws.Cells.Item(74, 11).Value2 = 872.2069  # K74: 714.025 -> 872.2069
$ws.Cells.Item(74, 12).Value2 = 15504  # L74: 21548.4 -> 15504
$ws.Cells.Item(74, 13).Value2 = 1.793099999999981  # M74: 159.975 -> 1.793099999999981
$ws.Cells.Item(74, 14).Value2 = -17252  # N74: -23296.4 -> -17252
$ws.Cells.Item(77, 8).Value2 = 3717.2778  # H77: 3028.9556 -> 3717.2778
$ws.Cells.Item(77, 9).Value2 = 872.2069  # I77: 714.025 -> 872.2069
$ws.Cells.Item(77, 10).Value2 = 15504  # J77: 21548.4 -> 15504
$ws.Cells.Item(77, 11).Value2 = 4361.0345  # K77: 3570.125 -> 4361.0345
$ws.Cells.Item(77, 12).Value2 = 77520  # L77: 107742 -> 77520
$ws.Cells.Item(77, 13).Value2 = 6.965500000000247  # M77: 797.875 -> 6.965500000000247
$ws.Cells.Item(77, 14).Value2 = -86256  # N77: -116478 -> -86256
$ws.Cells.Item(110, 8).Value2 = 2193.7932  # H110: 2247.1428 -> 2193.7932
$ws.Cells.Item(110, 9).Value2 = 2488.2354  # I110: 2600 -> 2488.2354
$ws.Cells.Item(110, 11).Value2 = 2488.2354  # K110: 2600 -> 2488.2354
$ws.Cells.Item(110, 13).Value2 = -443.2354  # M110: -555 -> -443.2354
$ws.Cells.Item(132, 8).Value2 = 9990.625  # H132: 6997.657 -> 9990.625
$ws.Cells.Item(132, 9).Value2 = 7235.278  # I132: 4565.933 -> 7235.278
$ws.Cells.Item(132, 10).Value2 = 18256.666  # J132: 21588 -> 18256.666
$ws.Cells.Item(132, 11).Value2 = 21705.834  # K132: 13697.799 -> 21705.834
$ws.Cells.Item(132, 12).Value2 = 54769.99800000001  # L132: 64764 -> 54769.99800000001
$ws.Cells.Item(132, 13).Value2 = -19175.834  # M132: -11167.799 -> -19175.834
$ws.Cells.Item(132, 14).Value2 = -59829.99800000001  # N132: -69824 -> -59829.99800000001
$ws.Cells.Item(136, 8).Value2 = 1145.4286  # H136: 1149.7206 -> 1145.4286
$ws.Cells.Item(136, 9).Value2 = 853.2406999999999  # I136: 865.5848999999999 -> 853.2406999999999
$ws.Cells.Item(136, 10).Value2 = 2131.5625  # J136: 2153.6667 -> 2131.5625
$ws.Cells.Item(136, 11).Value2 = 2559.7221  # K136: 2596.7547 -> 2559.7221
$ws.Cells.Item(136, 12).Value2 = 6394.6875  # L136: 6461.000100000001 -> 6394.6875
$ws.Cells.Item(136, 13).Value2 = -9.722099999999955  # M136: -46.75469999999996 -> -9.722099999999955
$ws.Cells.Item(136, 14).Value2 = -11494.6875  # N136: -11561.0001 -> -11494.6875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(87, 8).Value2 = 0  # H87: 38677 -> 0
$ws.Cells.Item(87, 10).Value2 = 0  # J87: 38677 -> 0
$ws.Cells.Item(87, 12).Value2 = 0  # L87: 38677 -> 0
$ws.Cells.Item(87, 14).Value2 = $null  # N87: -41173 -> None
$ws.Cells.Item(90, 8).Value2 = 0  # H90: 38677 -> 0
$ws.Cells.Item(90, 10).Value2 = 0  # J90: 38677 -> 0
$ws.Cells.Item(90, 12).Value2 = 0  # L90: 116031 -> 0
$ws.Cells.Item(90, 14).Value2 = $null  # N90: -128511 -> None
$ws.Cells.Item(105, 8).Value2 = 1843.2106  # H105: 1916.5 -> 1843.2106
$ws.Cells.Item(105, 9).Value2 = 1359  # I105: 1396.6666 -> 1359
$ws.Cells.Item(105, 10).Value2 = 2381.2222  # J105: 2306.375 -> 2381.2222
$ws.Cells.Item(105, 11).Value2 = 1359  # K105: 1396.6666 -> 1359
$ws.Cells.Item(105, 12).Value2 = 2381.2222  # L105: 2306.375 -> 2381.2222
$ws.Cells.Item(105, 13).Value2 = 388  # M105: 350.3334 -> 388
$ws.Cells.Item(105, 14).Value2 = -5875.2222  # N105: -5800.375 -> -5875.2222
$ws.Cells.Item(134, 8).Value2 = 1092.3334  # H134: 1220.4 -> 1092.3334
$ws.Cells.Item(134, 9).Value2 = 947.7143  # I134: 1064.3529 -> 947.7143
$ws.Cells.Item(134, 11).Value2 = 2843.1429  # K134: 3193.0587 -> 2843.1429
$ws.Cells.Item(134, 13).Value2 = -308.1428999999998  # M134: -658.0587000000005 -> -308.1428999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value2 = 5724.607  # H134: 2838.508 -> 5724.607
$ws.Cells.Item(134, 9).Value2 = 6652.136  # I134: 3080.453 -> 6652.136
$ws.Cells.Item(134, 10).Value2 = 2323.6667  # J134: 1556.2 -> 2323.6667
$ws.Cells.Item(134, 11).Value2 = 19956.408  # K134: 9241.359 -> 19956.408
$ws.Cells.Item(134, 12).Value2 = 6971.000100000001  # L134: 4668.6 -> 6971.000100000001
$ws.Cells.Item(134, 13).Value2 = -17421.408  # M134: -6706.359 -> -17421.408
$ws.Cells.Item(134, 14).Value2 = -12041.0001  # N134: -9738.6 -> -12041.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(88, 8).Value2 = 5000  # H88: 4785.7144 -> 5000
$ws.Cells.Item(88, 10).Value2 = 5000  # J88: 4785.7144 -> 5000
$ws.Cells.Item(88, 12).Value2 = 15000  # L88: 14357.1432 -> 15000
$ws.Cells.Item(88, 14).Value2 = -15856  # N88: -15213.1432 -> -15856
$ws.Cells.Item(91, 8).Value2 = 5000  # H91: 4785.7144 -> 5000
$ws.Cells.Item(91, 10).Value2 = 5000  # J91: 4785.7144 -> 5000
$ws.Cells.Item(91, 12).Value2 = 15000  # L91: 14357.1432 -> 15000
$ws.Cells.Item(91, 14).Value2 = -17964  # N91: -17321.1432 -> -17964
$ws.Cells.Item(134, 8).Value2 = 30242.375  # H134: 31549.078 -> 30242.375
$ws.Cells.Item(134, 9).Value2 = 40094.31  # I134: 39995.344 -> 40094.31
$ws.Cells.Item(134, 10).Value2 = 4269.091  # J134: 4333.3335 -> 4269.091
$ws.Cells.Item(134, 11).Value2 = 120282.93  # K134: 119986.032 -> 120282.93
$ws.Cells.Item(134, 12).Value2 = 12807.273  # L134: 13000.0005 -> 12807.273
$ws.Cells.Item(134, 13).Value2 = -115212.93  # M134: -114916.032 -> -115212.93
$ws.Cells.Item(134, 14).Value2 = -22947.273  # N134: -23140.0005 -> -22947.273
$ws.Cells.Item(137, 8).Value2 = 3191987.2  # H137: 4642199 -> 3191987.2
$ws.Cells.Item(137, 9).Value2 = 62218.89  # I137: 66018.234 -> 62218.89
$ws.Cells.Item(137, 10).Value2 = 7215975  # J137: 20201212 -> 7215975
$ws.Cells.Item(137, 11).Value2 = 186656.67  # K137: 198054.702 -> 186656.67
$ws.Cells.Item(137, 12).Value2 = 21647925  # L137: 60603636 -> 21647925
$ws.Cells.Item(137, 13).Value2 = -181556.67  # M137: -192954.702 -> -181556.67
$ws.Cells.Item(137, 14).Value2 = -21658125  # N137: -60613836 -> -21658125
$ws.Cells.Item(139, 8).Value2 = 325246.16  # H139: 368689.62 -> 325246.16
$ws.Cells.Item(139, 9).Value2 = 441338.8  # I139: 580620.5 -> 441338.8
$ws.Cells.Item(139, 10).Value2 = 2766.5557  # J139: 2627.182 -> 2766.5557
$ws.Cells.Item(139, 11).Value2 = 1324016.4  # K139: 1741861.5 -> 1324016.4
$ws.Cells.Item(139, 12).Value2 = 8299.667099999999  # L139: 7881.545999999999 -> 8299.667099999999
$ws.Cells.Item(139, 13).Value2 = -1318876.4  # M139: -1736721.5 -> -1318876.4
$ws.Cells.Item(139, 14).Value2 = -18579.6671  # N139: -18161.546 -> -18579.6671
$ws.Cells.Item(140, 8).Value2 = 19275.477  # H140: 20512.65 -> 19275.477
$ws.Cells.Item(140, 9).Value2 = 27691  # I140: 35421.637 -> 27691
$ws.Cells.Item(140, 10).Value2 = 2444.4285  # J140: 2290.5557 -> 2444.4285
$ws.Cells.Item(140, 11).Value2 = 83073  # K140: 106264.911 -> 83073
$ws.Cells.Item(140, 12).Value2 = 7333.2855  # L140: 6871.6671 -> 7333.2855
$ws.Cells.Item(140, 13).Value2 = -77893  # M140: -101084.911 -> -77893
$ws.Cells.Item(140, 14).Value2 = -17693.2855  # N140: -17231.6671 -> -17693.2855

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value2 = 2596.7646  # H80: 2596.7058 -> 2596.7646
$ws.Cells.Item(80, 9).Value2 = 2001.875  # I80: 2102.5 -> 2001.875
$ws.Cells.Item(80, 10).Value2 = 3125.5557  # J80: 3036 -> 3125.5557
$ws.Cells.Item(80, 11).Value2 = 2001.875  # K80: 2102.5 -> 2001.875
$ws.Cells.Item(80, 12).Value2 = 3125.5557  # L80: 3036 -> 3125.5557
$ws.Cells.Item(80, 13).Value2 = -1003.875  # M80: -1104.5 -> -1003.875
$ws.Cells.Item(80, 14).Value2 = -5121.5557  # N80: -5032 -> -5121.5557
$ws.Cells.Item(83, 8).Value2 = 2596.7646  # H83: 2596.7058 -> 2596.7646
$ws.Cells.Item(83, 9).Value2 = 2001.875  # I83: 2102.5 -> 2001.875
$ws.Cells.Item(83, 10).Value2 = 3125.5557  # J83: 3036 -> 3125.5557
$ws.Cells.Item(83, 11).Value2 = 10009.375  # K83: 10512.5 -> 10009.375
$ws.Cells.Item(83, 12).Value2 = 15627.7785  # L83: 15180 -> 15627.7785
$ws.Cells.Item(83, 13).Value2 = -5017.375  # M83: -5520.5 -> -5017.375
$ws.Cells.Item(83, 14).Value2 = -25611.7785  # N83: -25164 -> -25611.7785
$ws.Cells.Item(132, 8).Value2 = 4637.41  # H132: 3245.5173 -> 4637.41
$ws.Cells.Item(132, 9).Value2 = 6007.6924  # I132: 3740.682 -> 6007.6924
$ws.Cells.Item(132, 10).Value2 = 1896.8462  # J132: 1689.2858 -> 1896.8462
$ws.Cells.Item(132, 11).Value2 = 18023.0772  # K132: 11222.046 -> 18023.0772
$ws.Cells.Item(132, 12).Value2 = 5690.5386  # L132: 5067.857400000001 -> 5690.5386
$ws.Cells.Item(132, 13).Value2 = -15493.0772  # M132: -8692.045999999998 -> -15493.0772
$ws.Cells.Item(132, 14).Value2 = -10750.5386  # N132: -10127.8574 -> -10750.5386

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 611.48  # H22: 9280.083000000001 -> 611.48
$ws.Cells.Item(22, 9).Value2 = 407.64285  # I22: 880 -> 407.64285
$ws.Cells.Item(22, 10).Value2 = 870.9091  # J22: 15280.143 -> 870.9091
$ws.Cells.Item(22, 11).Value2 = 407.64285  # K22: 880 -> 407.64285
$ws.Cells.Item(22, 12).Value2 = 870.9091  # L22: 15280.143 -> 870.9091
$ws.Cells.Item(22, 13).Value2 = -112.64285  # M22: -585 -> -112.64285
$ws.Cells.Item(22, 14).Value2 = -1460.9091  # N22: -15870.143 -> -1460.9091
$ws.Cells.Item(27, 8).Value2 = 611.48  # H27: 9280.083000000001 -> 611.48
$ws.Cells.Item(27, 9).Value2 = 407.64285  # I27: 880 -> 407.64285
$ws.Cells.Item(27, 10).Value2 = 870.9091  # J27: 15280.143 -> 870.9091
$ws.Cells.Item(27, 11).Value2 = 407.64285  # K27: 880 -> 407.64285
$ws.Cells.Item(27, 12).Value2 = 870.9091  # L27: 15280.143 -> 870.9091
$ws.Cells.Item(27, 13).Value2 = -300.64285  # M27: -773 -> -300.64285
$ws.Cells.Item(27, 14).Value2 = -1084.9091  # N27: -15494.143 -> -1084.9091
$ws.Cells.Item(68, 8).Value2 = 1798.3103  # H68: 1850.5 -> 1798.3103
$ws.Cells.Item(68, 9).Value2 = 1807.55  # I68: 1885.3846 -> 1807.55
$ws.Cells.Item(68, 10).Value2 = 1777.7778  # J68: 1785.7142 -> 1777.7778
$ws.Cells.Item(68, 11).Value2 = 1807.55  # K68: 1885.3846 -> 1807.55
$ws.Cells.Item(68, 12).Value2 = 1777.7778  # L68: 1785.7142 -> 1777.7778
$ws.Cells.Item(68, 13).Value2 = -1058.55  # M68: -1136.3846 -> -1058.55
$ws.Cells.Item(68, 14).Value2 = -3275.7778  # N68: -3283.7142 -> -3275.7778
$ws.Cells.Item(71, 8).Value2 = 1798.3103  # H71: 1850.5 -> 1798.3103
$ws.Cells.Item(71, 9).Value2 = 1807.55  # I71: 1885.3846 -> 1807.55
$ws.Cells.Item(71, 10).Value2 = 1777.7778  # J71: 1785.7142 -> 1777.7778
$ws.Cells.Item(71, 11).Value2 = 9037.75  # K71: 9426.923000000001 -> 9037.75
$ws.Cells.Item(71, 12).Value2 = 8888.889000000001  # L71: 8928.571 -> 8888.889000000001
$ws.Cells.Item(71, 13).Value2 = -5293.75  # M71: -5682.923000000001 -> -5293.75
$ws.Cells.Item(71, 14).Value2 = -16376.889  # N71: -16416.571 -> -16376.889
$ws.Cells.Item(92, 8).Value2 = 32000  # H92: 0 -> 32000
$ws.Cells.Item(92, 10).Value2 = 32000  # J92: 0 -> 32000
$ws.Cells.Item(92, 12).Value2 = 32000  # L92: 0 -> 32000
$ws.Cells.Item(92, 14).Value2 = -36992  # N92: None -> -36992
$ws.Cells.Item(132, 8).Value2 = 4918.186  # H132: 5775.4165 -> 4918.186
$ws.Cells.Item(132, 9).Value2 = 6276.269  # I132: 8258.736999999999 -> 6276.269
$ws.Cells.Item(132, 10).Value2 = 2841.1177  # J132: 2999.9412 -> 2841.1177
$ws.Cells.Item(132, 11).Value2 = 18828.807  # K132: 24776.211 -> 18828.807
$ws.Cells.Item(132, 12).Value2 = 8523.3531  # L132: 8999.8236 -> 8523.3531
$ws.Cells.Item(132, 13).Value2 = -16298.807  # M132: -22246.211 -> -16298.807
$ws.Cells.Item(132, 14).Value2 = -13583.3531  # N132: -14059.8236 -> -13583.3531
$ws.Cells.Item(136, 8).Value2 = 3959.8096  # H136: 4149.55 -> 3959.8096
$ws.Cells.Item(136, 9).Value2 = 1398.3939  # I136: 1477.9678 -> 1398.3939
$ws.Cells.Item(136, 11).Value2 = 4195.1817  # K136: 4433.903399999999 -> 4195.1817
$ws.Cells.Item(136, 13).Value2 = -1645.1817  # M136: -1883.903399999999 -> -1645.1817

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value2 = 3364.4375  # H132: 2474.1365 -> 3364.4375
$ws.Cells.Item(132, 9).Value2 = 5778.391  # I132: 3536.9211 -> 5778.391
$ws.Cells.Item(132, 10).Value2 = 1143.6  # J132: 1031.7858 -> 1143.6
$ws.Cells.Item(132, 11).Value2 = 17335.173  # K132: 10610.7633 -> 17335.173
$ws.Cells.Item(132, 12).Value2 = 3430.8  # L132: 3095.3574 -> 3430.8
$ws.Cells.Item(132, 13).Value2 = -14805.173  # M132: -8080.763300000001 -> -14805.173
$ws.Cells.Item(132, 14).Value2 = -8490.799999999999  # N132: -8155.357400000001 -> -8490.799999999999
$ws.Cells.Item(135, 8).Value2 = 18733.334  # H135: 24268.111 -> 18733.334
$ws.Cells.Item(135, 9).Value2 = 10650  # I135: 15000 -> 10650
$ws.Cells.Item(135, 10).Value2 = 34900  # J135: 26916.143 -> 34900
$ws.Cells.Item(135, 11).Value2 = 10650  # K135: 15000 -> 10650
$ws.Cells.Item(135, 12).Value2 = 34900  # L135: 26916.143 -> 34900
$ws.Cells.Item(135, 13).Value2 = -5580  # M135: -9930 -> -5580
$ws.Cells.Item(135, 14).Value2 = -45040  # N135: -37056.143 -> -45040
$ws.Cells.Item(136, 8).Value2 = 950.0633  # H136: 981.0405 -> 950.0633
$ws.Cells.Item(136, 9).Value2 = 361.60376  # I136: 378.14 -> 361.60376
$ws.Cells.Item(136, 10).Value2 = 2149.6155  # J136: 2237.0833 -> 2149.6155
$ws.Cells.Item(136, 11).Value2 = 1084.81128  # K136: 1134.42 -> 1084.81128
$ws.Cells.Item(136, 12).Value2 = 6448.8465  # L136: 6711.249899999999 -> 6448.8465
$ws.Cells.Item(136, 13).Value2 = 1465.18872  # M136: 1415.58 -> 1465.18872
$ws.Cells.Item(136, 14).Value2 = -11548.8465  # N136: -11811.2499 -> -11548.8465

